$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.017, 0.553, 0.605, 53, 45, -91),
    @(-0.336, 0.215, 0.606, 122, -57, 30),
    @(-0.449, -0.105, -0.089, -132, -6, -150),
    @(-0.019, -0.188, -0.08599999999999999, 52, -68, -89),
    @(-0.176, -0.007, 0.695, 171, 27, -48),
    @(0.176, 0.051, 0.797, -1, 1, -51),
    @(0.08400000000000001, -0.296, 0.446, 141, 17, 113),
    @(0.429, 0.131, -0.08400000000000001, -30, 57, -149),
    @(0.046, 0.064, 0.315, -55, 42, 84),
    @(0.332, -0.152, -0.068, 93, -5, -143)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
